# Add season-record columns (Wins / Losses / Ties) to the BOS_2008 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (AD1:AF1) -----------------------------------------------------
# Copy the formatting already used by the other header cells (bold, border,
# centered) from A1 onto the three new header cells, then set their text.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Season record values (rows 2-49) --------------------------------------
# Every player row on this sheet belongs to the same team/season, so the
# team's overall win-loss-tie record (95-67-0) is repeated down each column.
$ws.Range("AD2:AD49").Value = 95
$ws.Range("AE2:AE49").Value = 67
$ws.Range("AF2:AF49").Value = 0

Write-Output "season record columns added"
